$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.229.52"
$ws.Range("E2").Value = "  -2.13%  "
$ws.Range("D3").Value = "2.881.15"
$ws.Range("E3").Value = "  -3.13%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "522.49"
$ws.Range("E5").Value = "  -3.50%  "
$ws.Range("D6").Value = "140.18"
$ws.Range("E6").Value = "  -7.58%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "0.546"
$ws.Range("E8").Value = "  -4.57%  "
$ws.Range("D9").Value = "2.884.42"
$ws.Range("E9").Value = "  -3.22%  "
$ws.Range("D10").Value = "0.106"
$ws.Range("E10").Value = "  -6.53%  "
$ws.Range("D11").Value = "5.96"
$ws.Range("E11").Value = "  -2.74%  "
$ws.Range("D12").Value = "0.356"
$ws.Range("E12").Value = "  -3.16%  "
$ws.Range("D13").Value = "3.387.62"
$ws.Range("E13").Value = "  -2.87%  "
$ws.Range("E14").Value = "  +2.06%  "
$ws.Range("D15").Value = "60.328.99"
$ws.Range("E15").Value = "  -2.08%  "
$ws.Range("D16").Value = "22.42"
$ws.Range("E16").Value = "  -5.30%  "
$ws.Range("D17").Value = "2.889.69"
$ws.Range("E17").Value = "  -2.99%  "
$ws.Range("D18").Value = "0.0000139"
$ws.Range("E18").Value = "  -5.23%  "
$ws.Range("D19").Value = "4.95"
$ws.Range("E19").Value = "  -4.30%  "
$ws.Range("D20").Value = "11.54"
$ws.Range("E20").Value = "  -4.13%  "
$ws.Range("D21").Value = "352.25"
$ws.Range("E21").Value = "  -7.48%  "
$ws.Range("D22").Value = "6.55"
$ws.Range("E22").Value = "  -1.76%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("D24").Value = "5.69"
$ws.Range("E24").Value = "  +0.56%  "
$ws.Range("D25").Value = "64.36"
$ws.Range("E25").Value = "  -1.91%  "
$ws.Range("D26").Value = "0.449"
$ws.Range("E26").Value = "  -4.77%  "
$ws.Range("D27").Value = "0.177"
$ws.Range("E27").Value = "  -6.32%  "
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").Value = "7.75"
$ws.Range("E29").Value = "  -5.98%  "
$ws.Range("D30").Value = "0.0₃0820"
$ws.Range("E30").Value = "  -12.92%  "
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").Value = "1.66"
$ws.Range("E32").Value = "  -3.19%  "
$ws.Range("D33").Value = "19.49"
$ws.Range("E33").Value = "  -4.98%  "
$ws.Range("D34").Value = "149.84"
$ws.Range("E34").Value = "  -7.04%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").Value = "5.52"
$ws.Range("E35").Value = "  -6.93%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "4.27"
$ws.Range("E36").Value = "  -8.62%  "
$ws.Range("D37").Value = "0.987"
$ws.Range("E37").Value = "  -8.43%  "
$ws.Range("E38").Value = "  -6.08%  "
$ws.Range("D39").Value = "37.52"
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("E40").Value = "  -6.51%  "
$ws.Range("D41").Value = "3.68"
$ws.Range("E41").Value = "  -6.17%  "
$ws.Range("D42").Value = "2.276.10"
$ws.Range("E42").Value = "  -5.78%  "
$ws.Range("D43").Value = "0.643"
$ws.Range("E43").Value = "  -3.86%  "
$ws.Range("D44").Value = "0.0577"
$ws.Range("E44").Value = "  -2.77%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "20.01"
$ws.Range("E46").Value = "  -9.96%  "
$ws.Range("D47").Value = "4.87"
$ws.Range("E47").Value = "  -5.10%  "
$ws.Range("D48").Value = "0.0235"
$ws.Range("E48").Value = "  -4.51%  "
$ws.Range("D49").Value = "10.32"
$ws.Range("E49").Value = "  -1.30%  "
$ws.Range("D50").Value = "0.0913"
$ws.Range("E50").Value = "  -4.30%  "
$ws.Range("D51").Value = "245.44"
$ws.Range("E51").Value = "  -7.88%  "
